# Clean and sort out translators (sbml part done)
# Update "Regulator Name" (column A) values for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value  = "JUN"
$ws.Range("A6").Value  = "FOS_DD"

$ws.Range("A8").Value  = "FOXP3"
$ws.Range("A11").Value = "NFKAPPAB"
$ws.Range("A12").Value = "AP1"

$ws.Range("A22").Value = "NFAT"
$ws.Range("A23").Value = "NFKAPPAB"

$ws.Range("A28").Value = "CD132"
$ws.Range("A29").Value = "CD25"

$ws.Range("A31").Value = "IL2_EX"
$ws.Range("A32").Value = "IL2R"

$ws.Range("A52").Value = "CD28"
$ws.Range("A53").Value = "PI3K"
$ws.Range("A55").Value = "IL2R"
$ws.Range("A56").Value = "TCR"

$ws.Range("A60").Value = "TCR"
$ws.Range("A61").Value = "MTORC2"

$ws.Range("A63").Value = "FOXP3"
$ws.Range("A64").Value = "PTEN"

$ws.Range("A67").Value = "CD28"
$ws.Range("A68").Value = "RAS"
$ws.Range("A69").Value = "IL2_EX"
$ws.Range("A70").Value = "IL2R"
